# Add a new row 38 of data to each of the four database sheets, mirroring
# the existing layout/format used by the prior rows (A = timestamp with the
# "YYYY-MM-DD HH:MM:SS" number format, B-E = hex text strings, F-I = numbers).

$wb = $excel.ActiveWorkbook

$dateFmt = "YYYY-MM-DD HH:MM:SS"
$newTimestamp = 45824.46601851852

# MID_LFT_#1 (sheet 1) - new row 38
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A38").Value = $newTimestamp
$ws1.Range("A38").NumberFormat = $dateFmt
$ws1.Range("B38").Value = "0x01,0x90"
$ws1.Range("C38").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws1.Range("D38").Value = "0x01,0x7C"
$ws1.Range("E38").Value = "0x07"
$ws1.Range("F38").Value = 400
$ws1.Range("G38").Value = 568631262647113000000000.0
$ws1.Range("H38").Value = 380
$ws1.Range("I38").Value = 7

# MID_LFT_#2 (sheet 2) - new row 38
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A38").Value = $newTimestamp
$ws2.Range("A38").NumberFormat = $dateFmt
$ws2.Range("B38").Value = "0x01,0x7c"
$ws2.Range("C38").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Range("D38").Value = "0x01,0x70"
$ws2.Range("E38").Value = "0x19"
$ws2.Range("F38").Value = 380
$ws2.Range("G38").Value = 568432987514711000000000.0
$ws2.Range("H38").Value = 368
$ws2.Range("I38").Value = 25

# MID_PLT_#1 (sheet 3) - new row 38
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A38").Value = $newTimestamp
$ws3.Range("A38").NumberFormat = $dateFmt
$ws3.Range("B38").Value = "0x00,0x6e"
$ws3.Range("C38").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws3.Range("D38").Value = "0x00,0x6B"
$ws3.Range("E38").Value = "0x15"
$ws3.Range("F38").Value = 110
$ws3.Range("G38").Value = 568631262647113000000000.0
$ws3.Range("H38").Value = 107
$ws3.Range("I38").Value = 15

# MID_PLT_#2 (sheet 4) - new row 38
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A38").Value = $newTimestamp
$ws4.Range("A38").NumberFormat = $dateFmt
$ws4.Range("B38").Value = "0x00,0x82"
$ws4.Range("C38").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws4.Range("D38").Value = "0x00,0x7F"
$ws4.Range("E38").Value = "0x9"
$ws4.Range("F38").Value = 130
$ws4.Range("G38").Value = 568631262647113000000000.0
$ws4.Range("H38").Value = 127
$ws4.Range("I38").Value = 9
